$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Preserve the "final row" (thick bottom) border style that currently
#        lives on C14 by copying it onto C19 BEFORE row 14 is re-formatted
#        as a normal interior row.
$ws.Range("C14").Copy()
$ws.Range("C19").PasteSpecial(-4122)

# --- 2. Re-format B14:C18 (old last row + 4 new rows) as normal interior
#        rows, using row 13 as the format template.
$ws.Range("B13:C13").Copy()
$ws.Range("B14:C18").PasteSpecial(-4122)

# --- 3. B19 keeps the plain interior left-column style (column B only).
$ws.Range("B13").Copy()
$ws.Range("B19").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 4. Row heights for the (previously unformatted) new rows.
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 30
$ws.Rows.Item(18).RowHeight = 30
$ws.Rows.Item(19).RowHeight = 30

# --- 5. Fill in the new data (Task IDs 11-15 + descriptions).
$ws.Range("B15").Value = 11
$ws.Range("C15").Value = "1º functionality"
$ws.Range("B16").Value = 12
$ws.Range("C16").Value = "2º functionality"
$ws.Range("B17").Value = 13
$ws.Range("C17").Value = "use cases"
$ws.Range("B18").Value = 14
$ws.Range("C18").Value = "demo video"
$ws.Range("B19").Value = 15
$ws.Range("C19").Value = "final delivarable"

# --- 6. Selection matches the author's cursor position after the edit.
$ws.Range("E18").Select()
